$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
@("M2","Fasl","Fas","ECs",1,0.3333333333333333,0.098742,0.296226,0.4941333075334537,0.4941333075334537,3,1,18.11506266666667,54.345188,0.5475965948184532,0.5475965948184532,1.788717517832,16.098457660488,0.2705857165916988,0.2705857165916988),
@("M2","Fasl","Fas","FAPs",1,0.3333333333333333,0.098742,0.296226,0.4941333075334537,0.4941333075334537,3,1,8.621912666666665,25.865738,0.2606300681353106,0.2606300681353106,0.8513449005319998,7.662104104787999,0.1287859976103704,0.1287859976103704),
@("M2","Fasl","Fas","M2",1,0.3333333333333333,0.098742,0.296226,0.4941333075334537,0.4941333075334537,3,1,3.165745,9.497235,0.09569667044284827,0.09569667044284827,0.31259199279,2.81332793511,0.04728691228586351,0.04728691228586351),
@("M2","Fasl","Fas","sCs",1,0.3333333333333333,0.098742,0.296226,0.4941333075334537,0.4941333075334537,3,1,3.178315666666667,9.534947,0.0960766666033877,0.09607666660338772,0.313833245558,2.824499210022,0.04747468104552088,0.04747468104552088),
@("sCs","Fasl","Fas","ECs",3,1,0.1010866666666667,0.30326,0.5058666924665464,0.5058666924665465,3,1,18.11506266666667,54.345188,0.5475965948184532,0.5475965948184532,1.831191301431111,16.48072171288,0.2770108782267545,0.2770108782267545),
@("sCs","Fasl","Fas","FAPs",3,1,0.1010866666666667,0.30326,0.5058666924665464,0.5058666924665465,3,1,8.621912666666665,25.865738,0.2606300681353106,0.2606300681353106,0.8715604117644443,7.84404370588,0.1318440705249402,0.1318440705249402),
@("sCs","Fasl","Fas","M2",3,1,0.1010866666666667,0.30326,0.5058666924665464,0.5058666924665465,3,1,3.165745,9.497235,0.09569667044284827,0.09569667044284827,0.3200146095666667,2.8801314861,0.04840975815698476,0.04840975815698478),
@("sCs","Fasl","Fas","sCs",3,1,0.1010866666666667,0.30326,0.5058666924665464,0.5058666924665465,3,1,3.178315666666667,9.534947,0.0960766666033877,0.09607666660338772,0.3212853363577778,2.89156802722,0.04860198555786684,0.04860198555786685)
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = 2 + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range($cols[$j] + $r).Value = $rowVals[$j]
    }
}
